$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This status string is shared by every row currently in "Ready for handoff"
# state, across the Overview summary sheet and each per-locale detail sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width change ---
# The zh-cn/de-de status columns on the Overview sheet, and the Status
# column on each locale sheet, narrow now that the status text is shorter
# ("Ready for handoff" -> "In Translation"), matching the regenerated
# report's column width of ~13.41 characters. Excel snaps ColumnWidth to
# its internal pixel grid, so this input is the value that lands on the
# closest on-grid width to the report's target.
$newStatusColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
